# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / LevePrice / LeveProfit values
# to rows across the ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 233
$ws.Range("I6").Value = 233
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 699
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -587
$ws.Range("N6").ClearContents()
$ws.Range("H94").Value = 11114778
$ws.Range("I94").Value = 12503014
$ws.Range("K94").Value = 12503014
$ws.Range("M94").Value = -12502563
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 461774.6
$ws.Range("I45").Value = 844230.25
$ws.Range("J45").Value = 2827.8
$ws.Range("K45").Value = 844230.25
$ws.Range("L45").Value = 2827.8
$ws.Range("M45").Value = -843853.25
$ws.Range("N45").Value = -3581.8
$ws.Range("H74").Value = 2165.923
$ws.Range("I74").Value = 1929.8334
$ws.Range("K74").Value = 1929.8334
$ws.Range("M74").Value = -1055.8334
$ws.Range("H77").Value = 2165.923
$ws.Range("I77").Value = 1929.8334
$ws.Range("K77").Value = 9649.166999999999
$ws.Range("M77").Value = -5281.166999999999
$ws.Range("H97").Value = 11496159
$ws.Range("I97").Value = 13890159
$ws.Range("K97").Value = 13890159
$ws.Range("M97").Value = -13889663
$ws.Range("H122").Value = 1747.6296
$ws.Range("I122").Value = 1617.95
$ws.Range("J122").Value = 2118.1428
$ws.Range("K122").Value = 4853.85
$ws.Range("L122").Value = 6354.428400000001
$ws.Range("M122").Value = -2403.85
$ws.Range("N122").Value = -11254.4284
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H129").Value = 100575
$ws.Range("J129").Value = 100575
$ws.Range("L129").Value = 100575
$ws.Range("N129").Value = -110575
$ws.Range("H132").Value = 1498.6904
$ws.Range("I132").Value = 1105.8
$ws.Range("K132").Value = 3317.4
$ws.Range("M132").Value = -787.3999999999996
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 30000
$ws.Range("J51").Value = 30000
$ws.Range("L51").Value = 30000
$ws.Range("N51").Value = -30982
$ws.Range("H86").Value = 2999
$ws.Range("I86").Value = 1628.8572
$ws.Range("J86").Value = 4597.5
$ws.Range("K86").Value = 1628.8572
$ws.Range("L86").Value = 4597.5
$ws.Range("M86").Value = -505.8571999999999
$ws.Range("N86").Value = -6843.5
$ws.Range("H89").Value = 2999
$ws.Range("I89").Value = 1628.8572
$ws.Range("J89").Value = 4597.5
$ws.Range("K89").Value = 8144.286
$ws.Range("L89").Value = 22987.5
$ws.Range("M89").Value = -2528.286
$ws.Range("N89").Value = -34219.5
$ws.Range("H99").Value = 3777.6667
$ws.Range("I99").Value = 4033.1667
$ws.Range("J99").Value = 3266.6667
$ws.Range("K99").Value = 4033.1667
$ws.Range("L99").Value = 3266.6667
$ws.Range("M99").Value = -2535.1667
$ws.Range("N99").Value = -6262.6667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 15000
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15470
$ws.Range("H32").Value = 4692.375
$ws.Range("I32").Value = 1791.2858
$ws.Range("K32").Value = 1791.2858
$ws.Range("M32").Value = -1475.2858
$ws.Range("H58").Value = 1273.9231
$ws.Range("I58").Value = 1069.4
$ws.Range("J58").Value = 1955.6666
$ws.Range("K58").Value = 1069.4
$ws.Range("L58").Value = 1955.6666
$ws.Range("M58").Value = -866.4000000000001
$ws.Range("N58").Value = -2361.6666
$ws.Range("H136").Value = 1273.9231
$ws.Range("I136").Value = 1069.4
$ws.Range("J136").Value = 1955.6666
$ws.Range("K136").Value = 3208.2
$ws.Range("L136").Value = 5866.9998
$ws.Range("M136").Value = -658.2000000000003
$ws.Range("N136").Value = -10966.9998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 5000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18733.938
$ws.Range("I7").Value = 30012.555
$ws.Range("J7").Value = 4232.857
$ws.Range("K7").Value = 30012.555
$ws.Range("L7").Value = 4232.857
$ws.Range("M7").Value = -29900.555
$ws.Range("N7").Value = -4456.857
$ws.Range("H22").Value = 1295.375
$ws.Range("I22").Value = 904.25
$ws.Range("K22").Value = 904.25
$ws.Range("M22").Value = -609.25
$ws.Range("H27").Value = 1295.375
$ws.Range("I27").Value = 904.25
$ws.Range("K27").Value = 904.25
$ws.Range("M27").Value = -797.25
$ws.Range("H93").Value = 2638.1333
$ws.Range("I93").Value = 2614.3333
$ws.Range("J93").Value = 2673.8333
$ws.Range("K93").Value = 2614.3333
$ws.Range("L93").Value = 2673.8333
$ws.Range("M93").Value = -1366.3333
$ws.Range("N93").Value = -5169.8333
$ws.Range("H100").Value = 6076.269
$ws.Range("J100").Value = 4317
$ws.Range("L100").Value = 4317
$ws.Range("N100").Value = -5399
$ws.Range("H122").Value = 8212.315000000001
$ws.Range("I122").Value = 8631.412
$ws.Range("K122").Value = 25894.236
$ws.Range("M122").Value = -23444.236
$ws.Range("H126").Value = 18733.938
$ws.Range("I126").Value = 30012.555
$ws.Range("J126").Value = 4232.857
$ws.Range("K126").Value = 90037.66500000001
$ws.Range("L126").Value = 12698.571
$ws.Range("M126").Value = -87567.66500000001
$ws.Range("N126").Value = -17638.571
$ws.Range("H132").Value = 336028.2
$ws.Range("I132").Value = 501835.56
$ws.Range("J132").Value = 4413.4
$ws.Range("K132").Value = 1505506.68
$ws.Range("L132").Value = 13240.2
$ws.Range("M132").Value = -1502976.68
$ws.Range("N132").Value = -18300.2
$ws.Range("H136").Value = 3740.5
$ws.Range("I136").Value = 2598
$ws.Range("K136").Value = 7794
$ws.Range("M136").Value = -5244
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6457.143
$ws.Range("I81").Value = 6457.143
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 12914.286
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -11853.286
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 6457.143
$ws.Range("I84").Value = 6457.143
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 64571.43
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -59267.43
$ws.Range("N84").ClearContents()
$ws.Range("H96").Value = 4717.6
$ws.Range("J96").Value = 3393.75
$ws.Range("L96").Value = 3393.75
$ws.Range("N96").Value = -6139.75
$ws.Range("H132").Value = 27675.324
$ws.Range("I132").Value = 26114.686
$ws.Range("J132").Value = 38599.8
$ws.Range("K132").Value = 78344.058
$ws.Range("L132").Value = 115799.4
$ws.Range("M132").Value = -75814.058
$ws.Range("N132").Value = -120859.4
